$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166318655014038
$ws.Range("B1").Value = 2.435633420944214
$ws.Range("D1").Value = 2.368980169296265
$ws.Range("E1").Value = 1.234718441963196
